# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the daily conversion note text (cell A1) ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$oldText = $ws1.Range("A1").Value2
$newText = $oldText.Replace(
    "1000 Bs = 2.68 = 10080.43 pesos",
    "1000 Bs = 2.54 = 9536.3 pesos"
).Replace(
    "10080.43 pesos = 2.68 = 951.74 Bs",
    "9536.3 pesos = 2.52 = 930.98 Bs"
)
$ws1.Range("A1").Value = $newText

# --- tasas: update N10/O10 and N12/O12 ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 393.9
$ws2.Range("O10").Value = 3756.35

$ws2.Range("N12").Value = 3780
$ws2.Range("O12").Value = 369.02
